# Add a new "mh-toc-parity" resource row to the Observations summary sheet.
#
# Columns (row 1 headers): Profile | Name | Category Code | Category VS |
#   Code | Code VS | Time Types | Value Types | Data Absent Reason |
#   Body Site | Method

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the header row's formatting down into row 2 first, so every new
# cell (including the ones that stay blank) picks up the same style (s="2")
# used across row 1.
$ws.Range("A1:K1").Copy()
$ws.Range("A2:K2").PasteSpecial(-4122)

$jHat = [char]0x135

$ws.Range("A2").Value = "mh-toc-parity"
$ws.Range("B2").Value = "MH TOC Parity Profile"
# Category Code (C2), Category VS (D2) -> left blank
$ws.Range("E2").Value = "LOINC#11977-6"
# Code VS (F2) -> left blank
$ws.Range("G2").Value = "dateTime$jHat, Period$jHat, Timing$jHat, instant$jHat"
$ws.Range("H2").Value = "integer$jHat, time$jHat, dateTime$jHat, Period$jHat"
$ws.Range("I2").Value = "optional"
# Body Site (J2), Method (K2) -> left blank
